$p = $ppt.ActivePresentation

# Delete slide 3 ("Fejlesztés") - sldId 260 - which was removed from the deck.
$p.Slides.Item(3).Delete()
